$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 43794.432638888888
$ws.Range("B2").Value = 43794.53402777778
$ws.Range("C2").Value = 43796.884722222225

$ws.Range("A3").Value = 43794.511805555558
$ws.Range("B3").Value = 43794.520833333336
$ws.Range("C3").Value = 43799.333333333336

$ws.Range("A4").Value = 43792.583333333336
$ws.Range("B4").Value = 43793.65625
$ws.Range("C4").Value = 43797.916666666664

$ws.Range("E5").Select()
